$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2741.52
$ws.Range("I15").Value = 2741.52
$ws.Range("K15").Value = 8224.559999999999
$ws.Range("M15").Value = -8055.559999999999
$ws.Range("H41").Value = 2018.7222
$ws.Range("I41").Value = 1918.2307
$ws.Range("J41").Value = 2280
$ws.Range("K41").Value = 1918.2307
$ws.Range("L41").Value = 2280
$ws.Range("M41").Value = -1478.2307
$ws.Range("N41").Value = -3160
$ws.Range("H125").Value = 2298.8235
$ws.Range("I125").Value = 3466.6667
$ws.Range("K125").Value = 31200.0003
$ws.Range("M125").Value = -28740.0003
$ws.Range("H135").Value = 24390826
$ws.Range("I135").Value = 237.28572
$ws.Range("J135").Value = 166669250
$ws.Range("K135").Value = 2135.57148
$ws.Range("L135").Value = 1500023250
$ws.Range("M135").Value = 399.4285199999999
$ws.Range("N135").Value = -1500028320
$ws.Range("H137").Value = 1100.8611
$ws.Range("I137").Value = 810.6739
$ws.Range("J137").Value = 1614.2693
$ws.Range("K137").Value = 2432.0217
$ws.Range("L137").Value = 4842.8079
$ws.Range("M137").Value = 117.9782999999998
$ws.Range("N137").Value = -9942.8079
$ws.Range("H138").Value = 1322.69
$ws.Range("I138").Value = 785.42426
$ws.Range("J138").Value = 1587.3135
$ws.Range("K138").Value = 2356.27278
$ws.Range("L138").Value = 4761.9405
$ws.Range("M138").Value = 2783.72722
$ws.Range("N138").Value = -15041.9405
$ws.Range("H141").Value = 632.7308
$ws.Range("I141").Value = 552.125
$ws.Range("K141").Value = 1656.375
$ws.Range("M141").Value = 3523.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1073.5
$ws.Range("I2").Value = 711
$ws.Range("J2").Value = 1113.7778
$ws.Range("K2").Value = 711
$ws.Range("L2").Value = 1113.7778
$ws.Range("M2").Value = -598
$ws.Range("N2").Value = -1339.7778
$ws.Range("H32").Value = 3547.035
$ws.Range("I32").Value = 3223.4167
$ws.Range("J32").Value = 5211.357
$ws.Range("K32").Value = 3223.4167
$ws.Range("L32").Value = 5211.357
$ws.Range("M32").Value = -2936.4167
$ws.Range("N32").Value = -5785.357
$ws.Range("H61").Value = 22223356
$ws.Range("I61").Value = 32258830
$ws.Range("J61").Value = 1950.8572
$ws.Range("K61").Value = 32258830
$ws.Range("L61").Value = 1950.8572
$ws.Range("M61").Value = -32258618
$ws.Range("N61").Value = -2374.8572
$ws.Range("H74").Value = 1875
$ws.Range("I74").Value = 1320.5834
$ws.Range("J74").Value = 3205.6
$ws.Range("K74").Value = 1320.5834
$ws.Range("L74").Value = 3205.6
$ws.Range("M74").Value = -446.5834
$ws.Range("N74").Value = -4953.6
$ws.Range("H77").Value = 1875
$ws.Range("I77").Value = 1320.5834
$ws.Range("J77").Value = 3205.6
$ws.Range("K77").Value = 6602.916999999999
$ws.Range("L77").Value = 16028
$ws.Range("M77").Value = -2234.916999999999
$ws.Range("N77").Value = -24764
$ws.Range("H116").Value = 1073.5
$ws.Range("I116").Value = 711
$ws.Range("J116").Value = 1113.7778
$ws.Range("K116").Value = 711
$ws.Range("L116").Value = 1113.7778
$ws.Range("M116").Value = 1583
$ws.Range("N116").Value = -5701.7778
$ws.Range("H122").Value = 1617.909
$ws.Range("I122").Value = 2079.5715
$ws.Range("J122").Value = 810
$ws.Range("K122").Value = 6238.7145
$ws.Range("L122").Value = 2430
$ws.Range("M122").Value = -3788.7145
$ws.Range("N122").Value = -7330
$ws.Range("H132").Value = 1851.2727
$ws.Range("I132").Value = 1645.25
$ws.Range("J132").Value = 2168.2307
$ws.Range("K132").Value = 4935.75
$ws.Range("L132").Value = 6504.6921
$ws.Range("M132").Value = -2405.75
$ws.Range("N132").Value = -11564.6921
$ws.Range("H136").Value = 22223356
$ws.Range("I136").Value = 32258830
$ws.Range("J136").Value = 1950.8572
$ws.Range("K136").Value = 96776490
$ws.Range("L136").Value = 5852.571599999999
$ws.Range("M136").Value = -96773940
$ws.Range("N136").Value = -10952.5716
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1073.5
$ws.Range("I3").Value = 711
$ws.Range("J3").Value = 1113.7778
$ws.Range("K3").Value = 711
$ws.Range("L3").Value = 1113.7778
$ws.Range("M3").Value = -597
$ws.Range("N3").Value = -1341.7778
$ws.Range("H134").Value = 4143.108
$ws.Range("I134").Value = 1008.67645
$ws.Range("K134").Value = 3026.02935
$ws.Range("M134").Value = -491.0293500000002
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1647.4807
$ws.Range("I31").Value = 1481.9783
$ws.Range("J31").Value = 2916.3333
$ws.Range("K31").Value = 1481.9783
$ws.Range("L31").Value = 2916.3333
$ws.Range("M31").Value = -1186.9783
$ws.Range("N31").Value = -3506.3333
$ws.Range("H34").Value = 1647.4807
$ws.Range("I34").Value = 1481.9783
$ws.Range("J34").Value = 2916.3333
$ws.Range("K34").Value = 1481.9783
$ws.Range("L34").Value = 2916.3333
$ws.Range("M34").Value = -1279.9783
$ws.Range("N34").Value = -3320.3333
$ws.Range("H58").Value = 852.30554
$ws.Range("I58").Value = 764.6774
$ws.Range("J58").Value = 1395.6
$ws.Range("K58").Value = 764.6774
$ws.Range("L58").Value = 1395.6
$ws.Range("M58").Value = -561.6774
$ws.Range("N58").Value = -1801.6
$ws.Range("H132").Value = 4544.1284
$ws.Range("I132").Value = 5402.893
$ws.Range("K132").Value = 16208.679
$ws.Range("M132").Value = -13678.679
$ws.Range("H134").Value = 10639594
$ws.Range("I134").Value = 1445.2188
$ws.Range("J134").Value = 33334310
$ws.Range("K134").Value = 4335.6564
$ws.Range("L134").Value = 100002930
$ws.Range("M134").Value = -1800.6564
$ws.Range("N134").Value = -100008000
$ws.Range("H136").Value = 852.30554
$ws.Range("I136").Value = 764.6774
$ws.Range("J136").Value = 1395.6
$ws.Range("K136").Value = 2294.0322
$ws.Range("L136").Value = 4186.799999999999
$ws.Range("M136").Value = 255.9677999999999
$ws.Range("N136").Value = -9286.799999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 643.075
$ws.Range("I113").Value = 468.6
$ws.Range("J113").Value = 668
$ws.Range("K113").Value = 1405.8
$ws.Range("L113").Value = 2004
$ws.Range("M113").Value = 764.1999999999998
$ws.Range("N113").Value = -6344
$ws.Range("H131").Value = 14286831
$ws.Range("J131").Value = 1209.6167
$ws.Range("L131").Value = 3628.8501
$ws.Range("N131").Value = -13708.8501
$ws.Range("H137").Value = 20838774
$ws.Range("I137").Value = 46876784
$ws.Range("J137").Value = 8367.15
$ws.Range("K137").Value = 140630352
$ws.Range("L137").Value = 25101.45
$ws.Range("M137").Value = -140625252
$ws.Range("N137").Value = -35301.45
$ws.Range("H140").Value = 19659.932
$ws.Range("I140").Value = 52477.4
$ws.Range("J140").Value = 2830.4614
$ws.Range("K140").Value = 157432.2
$ws.Range("L140").Value = 8491.3842
$ws.Range("M140").Value = -152252.2
$ws.Range("N140").Value = -18851.3842
$ws.Range("H141").Value = 1715.2667
$ws.Range("I141").Value = 1715.2667
$ws.Range("K141").Value = 5145.800099999999
$ws.Range("M141").Value = 34.19990000000053
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2374.75
$ws.Range("I126").Value = 1712.1333
$ws.Range("K126").Value = 5136.3999
$ws.Range("M126").Value = -2666.3999
$ws.Range("H132").Value = 3157.68
$ws.Range("I132").Value = 2943.7646
$ws.Range("J132").Value = 3612.25
$ws.Range("K132").Value = 8831.293799999999
$ws.Range("L132").Value = 10836.75
$ws.Range("M132").Value = -6301.293799999999
$ws.Range("N132").Value = -15896.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 21978.3
$ws.Range("I132").Value = 1555.8518
$ws.Range("J132").Value = 45952.477
$ws.Range("K132").Value = 4667.555399999999
$ws.Range("L132").Value = 137857.431
$ws.Range("M132").Value = -2137.555399999999
$ws.Range("N132").Value = -142917.431
$ws.Range("H136").Value = 4164.5137
$ws.Range("I136").Value = 4906.6787
$ws.Range("J136").Value = 1855.5555
$ws.Range("K136").Value = 14720.0361
$ws.Range("L136").Value = 5566.666499999999
$ws.Range("M136").Value = -12170.0361
$ws.Range("N136").Value = -10666.6665
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1186.6604
$ws.Range("I132").Value = 1076.1786
$ws.Range("J132").Value = 1310.4
$ws.Range("K132").Value = 3228.5358
$ws.Range("L132").Value = 3931.2
$ws.Range("M132").Value = -698.5357999999997
$ws.Range("N132").Value = -8991.200000000001
$ws.Range("H136").Value = 414.60974
$ws.Range("I136").Value = 343.35483
$ws.Range("K136").Value = 1030.06449
$ws.Range("M136").Value = 1519.93551
$ws.Range("N136").Value = -7006.5

Write-Host "Done applying changes."